$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D ("Price") hold plain text values in the source workbook
# (e.g. "68.273.77", "598.43"). Because some of the new values look like
# pure numbers (e.g. "598.43"), a plain .Value assignment would make Excel
# auto-convert them to numeric cells. To keep them as text (matching the
# original cell type), each cell is briefly switched to the Text ("@")
# number format before the assignment and then restored to the default
# "Normal" style afterwards so no extra formatting is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.321.33'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.21%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.652.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.20%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.42%  '

$ws.Range('E8').Value = '  -0.30%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.651.91'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.16%  '

$ws.Range('E10').Value = '  +7.38%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.158'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.01%  '

$ws.Range('E12').Value = '  +0.45%  '

$ws.Range('E13').Value = '  +1.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.02%  '

$ws.Range('E15').Value = '  +2.37%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.132.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.34%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.293.97'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.02%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.647.79'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.68%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.35%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '364.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.12%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.92%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.44%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.42%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.91%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '75.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.34%  '

$ws.Range('E26').Value = '  +0.01%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.91%  '

$ws.Range('E28').Value = '  +1.38%  '

$ws.Range('E29').Value = '  +1.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '560.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.26%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.17%  '

$ws.Range('E33').Value = '  +0.95%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.130'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.61%  '

$ws.Range('E36').Value = '  -0.07%  '

$ws.Range('E37').Value = '  +2.80%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '161.25'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.81%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.374'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.75%  '

$ws.Range('E41').Value = '  -0.87%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.38%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₆0341'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.11%  '

$ws.Range('E44').Value = '  -0.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.89%  '

$ws.Range('E46').Value = '  +0.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.35%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '159.28'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.24%  '

$ws.Range('E49').Value = '  +1.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.06'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.21%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.70'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.22%  '
